$d = $word.ActiveDocument

function Find-ParagraphIndexByText($doc, $needle) {
    $paras = $doc.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $t = $paras.Item($i).Range.Text
        if ($t -like $needle) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------
# Change 1: insert a new paragraph "Consignes données : " just
# before "Questions type pour le ressenti global de la partie :"
# ---------------------------------------------------------------
$idx = Find-ParagraphIndexByText $d "*ressenti global de la partie*"
$targetPara = $d.Paragraphs.Item($idx)
$targetPara.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs.Item($idx)
$newRange = $newPara.Range

$fragment = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Corpsdetexte"/><w:rPr><w:rFonts w:ascii="System Font" w:hAnsi="System Font" w:cs="System Font"/><w:bCs/><w:color w:val="FF0000"/><w:kern w:val="0"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:bCs/></w:rPr><w:t>Consignes données</w:t></w:r><w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space="preserve"> : </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newRange.InsertXML($fragment)

# ---------------------------------------------------------------
# Change 2: remove the stray <w:lastRenderedPageBreak/> in front of
# "Questions type pour le visionnage de la vidéo :"
# ---------------------------------------------------------------
$idx2 = Find-ParagraphIndexByText $d "*visionnage de la vid*"
$para2 = $d.Paragraphs.Item($idx2)
$range2 = $para2.Range
$clean2 = $range2.WordOpenXML
$range2.InsertXML($clean2)

# ---------------------------------------------------------------
# Change 3: collapse the two consecutive empty paragraphs after the
# comment reference (id 1) down to a single empty paragraph
# ---------------------------------------------------------------
$idx3 = Find-ParagraphIndexByText $d "*entretien et poser nos questions*"
$emptyPara = $d.Paragraphs.Item($idx3 + 2)
$emptyPara.Range.Delete()

Write-Host "done"
